# Refresh scraped cryptocurrency price/volume figures (GitHub Actions run 2024-06-23)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A few price cells round-trip through Excel's automatic number detection and would
# lose a trailing zero (e.g. "578.00" -> 578) if assigned as plain text, so those
# specific cells are pinned to the Text number format first, just like a user would
# do by pre-formatting the cell before typing a value that should stay literal.
foreach ($r in @(5, 27, 39, 46, 48)) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

$ws.Range("D2").Value = "63.606.38"
$ws.Range("E2").Value = "  -1.10%  "
$ws.Range("D3").Value = "3.427.11"
$ws.Range("E3").Value = "  -2.12%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "578.00"
$ws.Range("E5").Value = "  -2.18%  "
$ws.Range("D6").Value = "129.21"
$ws.Range("E6").Value = "  -3.56%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -1.40%  "
$ws.Range("D9").Value = "7.56"
$ws.Range("E9").Value = "  +2.81%  "
$ws.Range("E10").Value = "  -1.19%  "
$ws.Range("E11").Value = "  -1.21%  "
$ws.Range("D12").Value = "4.008.42"
$ws.Range("E12").Value = "  -2.22%  "
$ws.Range("E13").Value = "  -0.34%  "
$ws.Range("E14").Value = "  -2.78%  "
$ws.Range("D15").Value = "3.426.58"
$ws.Range("E15").Value = "  -2.15%  "
$ws.Range("D16").Value = "63.594.91"
$ws.Range("E16").Value = "  -1.26%  "
$ws.Range("D17").Value = "25.12"
$ws.Range("E17").Value = "  -2.04%  "
$ws.Range("D18").Value = "9.83"
$ws.Range("E18").Value = "  -0.39%  "
$ws.Range("E19").Value = "  -2.12%  "
$ws.Range("E20").Value = "  -1.90%  "
$ws.Range("D21").Value = "383.91"
$ws.Range("E21").Value = "  -2.29%  "
$ws.Range("E22").Value = "  -1.75%  "
$ws.Range("D23").Value = "3.559.78"
$ws.Range("E23").Value = "  -2.27%  "
$ws.Range("D24").Value = "73.82"
$ws.Range("E24").Value = "  -1.09%  "
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("E26").Value = "  -5.16%  "
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.31%  "
$ws.Range("D28").Value = "2.18"
$ws.Range("E28").Value = "  -3.33%  "
$ws.Range("D29").Value = "7.02"
$ws.Range("E29").Value = "  -4.74%  "
$ws.Range("D30").Value = "7.89"
$ws.Range("E30").Value = "  -4.02%  "
$ws.Range("E31").Value = "  -0.45%  "
$ws.Range("E32").Value = "  -4.49%  "
$ws.Range("D33").Value = "3.452.38"
$ws.Range("E33").Value = "  -2.00%  "
$ws.Range("D35").Value = "22.67"
$ws.Range("E35").Value = "  -3.37%  "
$ws.Range("D36").Value = "5.18"
$ws.Range("E36").Value = "  +1.01%  "
$ws.Range("D37").Value = "6.74"
$ws.Range("E37").Value = "  -1.92%  "
$ws.Range("D38").Value = "164.19"
$ws.Range("E38").Value = "  -1.87%  "
$ws.Range("D39").Value = "1.50"
$ws.Range("E39").Value = "  -2.75%  "
$ws.Range("D40").Value = "0.0763"
$ws.Range("E40").Value = "  -2.16%  "
$ws.Range("D41").Value = "0.784"
$ws.Range("E41").Value = "  -3.13%  "
$ws.Range("E42").Value = "  -0.08%  "
$ws.Range("D43").Value = "41.36"
$ws.Range("E43").Value = "  -1.04%  "
$ws.Range("D44").Value = "4.31"
$ws.Range("E44").Value = "  -1.91%  "
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").Value = "1.59"
$ws.Range("E45").Value = "  -4.09%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "23.20"
$ws.Range("E46").Value = "  -7.91%  "
$ws.Range("E47").Value = "  -5.99%  "
$ws.Range("D48").Value = "6.70"
$ws.Range("E48").Value = "  -0.77%  "
$ws.Range("D49").Value = "0.882"
$ws.Range("E49").Value = "  -1.12%  "
$ws.Range("D50").Value = "2.270.18"
$ws.Range("E50").Value = "  -4.57%  "
$ws.Range("D51").Value = "0.0251"
$ws.Range("E51").Value = "  -2.83%  "
